# Generate Report for Handoff
# The f5348948-...md file is now "Ready for handoff" with fresh handoff
# timestamps, and an Error Detail noting the handback file is stale.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$overviewDate = "2016-10-17 13:51:16"
$zhHandoffDate = "2016-10-17 13:50:54"
$deHandoffDate = "2016-10-17 13:51:16"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e861ea329665b9b5b0879684cabd4ecd9939d1e/e2e/f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7f58b95fd07e5bc7780ee4fd4e5b19ce267203cb/e2e/f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md."

# --- "Overview" sheet: row 3 is the f5348948-...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $overviewDate

# Excel's stored column <col width> = ColumnWidth + 5/MaximumDigitWidth
# (here 5/6), so back the requested width off by that padding to land on
# an exact stored width of 40.
$colWidth = 40 - (5 / 6)

# --- "zh-cn" sheet: row 3 is the f5348948-...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = $zhHandoffDate
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $colWidth

# --- "de-de" sheet: row 3 is the f5348948-...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = $deHandoffDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $colWidth
